$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "TC_02"
$ws.Range("B3").Value = 'Login->To Verify that Successfully landed user role-based landing page, when user Clicking on the "Adva pro Login" button. '
$ws.Range("F3").Value = "no"
